$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Completa el texto..." paragraph -> expanded wording, and the
#    Word "last edit" bookmark (_GoBack) now wraps this paragraph's
#    text (it used to sit in the "Tomado de" paragraph below).
# ------------------------------------------------------------------
$oldText1 = "Completa el texto, para ello arrastra las palabras al lugar adecuado."
$newText1 = "Ejercicio de completar el texto. Para ello, el estudiante debe arrastrar las palabras al lugar adecuado."

$paras = $d.Paragraphs
$target1 = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq $oldText1) {
        $target1 = $p
        break
    }
}

$r1 = $target1.Range
$r1.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $newText1, 1) | Out-Null

# Re-fetch the paragraph (content changed). The bookmark spans the run
# text plus the paragraph mark (matching Word's normal "last edit"
# range, which swallows the following pilcrow too).
$target1 = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq $newText1) {
        $target1 = $p
        break
    }
}
$bmRange = $d.Range($target1.Range.Start, $target1.Range.End + 1)

# Remove the old, now-stale _GoBack bookmark (wherever it is) before
# adding the new one back at the actual edit location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ------------------------------------------------------------------
# 2) "Tomado de" paragraph -> the two runs "Tomado de" + " " collapse
#    into a single run "Tomado de " once the bookmark that used to
#    separate them is gone.
# ------------------------------------------------------------------
$oldText2 = "Tomado de "
$paras = $d.Paragraphs
$target2 = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -like "Tomado de *") {
        $target2 = $p
        break
    }
}
$r2 = $target2.Range
$r2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $oldText2, 1) | Out-Null

Write-Output "done"
